$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "This movie was absolutely fantastic! I loved every minute of it."
